$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct text assignments (values that are unambiguously non-numeric, e.g.
# contain a "%" sign, multiple dots, or are plain names/URLs) -- safe to set
# via .Value directly since Excel will not misinterpret them as numbers.
$ws.Range("D2").Value = "29.055.99"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.834.91"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("E6").Value = "  -3.16%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "1.831.79"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("E17").Value = "  -3.07%  "
$ws.Range("D18").Value = "29.031.05"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "2.082.56"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.214.01"
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "1.980.89"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("E48").Value = "  -5.14%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("E51").Value = "  +0.67%  "

# Numeric-looking text values (e.g. "1.000", "242.45") must be written as literal
# text, not converted to numbers, to preserve formatting (trailing zeros etc.),
# matching the original inline-string cell content exactly.
# Writing them directly (even with a leading apostrophe) causes Excel/COM to
# stamp the cell with a distinct "quote prefix" style, which the source workbook
# does not have. Routing the text through a scratch cell + Copy/PasteSpecial
# (values only) carries over the text value without carrying over that style,
# keeping cell formatting identical to the original.
$scratch = $ws.Range("ZZ1")
$scratch.Value = "'1.000"
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$scratch.Value = "'242.45"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Value = "'0.6149"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Value = "'1.002"
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$scratch.Value = "'0.07477"
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$scratch.Value = "'0.2923"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$scratch.Value = "'0.07681"
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$scratch.Value = "'4.997"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$scratch.Value = "'0.6719"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$scratch.Value = "'82.63"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Value = "'0.000009155"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$scratch.Value = "'5.910"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$scratch.Value = "'232.17"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$scratch.Value = "'1.002"
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Value = "'7.189"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$scratch.Value = "'1.002"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$scratch.Value = "'158.71"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$scratch.Value = "'0.1396"
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$scratch.Value = "'8.498"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Value = "'17.81"
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$scratch.Value = "'1.498"
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$scratch.Value = "'4.151"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Value = "'4.116"
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$scratch.Value = "'0.05507"
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$scratch.Value = "'1.199"
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$scratch.Value = "'1.833"
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$scratch.Value = "'0.7368"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$scratch.Value = "'1.141"
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$scratch.Value = "'2.658"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Value = "'2.773"
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Value = "'0.01778"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Value = "'6.429"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Value = "'0.8890"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Value = "'101.90"
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$scratch.Value = "'0.5087"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$scratch.Value = "'0.00000000120"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$scratch.Value = "'0.4064"
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Value = "'9.139"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Value = "'0.05828"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = 0
